$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 10 (ALC)
$ws.Range("H10").Value = 19980
$ws.Range("J10").Value = 19980
$ws.Range("L10").Value = 19980
$ws.Range("N10").Value = -20566

# Row 17 (ALC)
$ws.Range("H17").Value = 1698.1072
$ws.Range("J17").Value = 1698.1072
$ws.Range("L17").Value = 5094.321599999999
$ws.Range("N17").Value = -5430.321599999999

# Row 28 (ALC)
$ws.Range("H28").Value = 3249.75
$ws.Range("I28").Value = 3249.75
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 3249.75
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -2764.75
$ws.Range("N28").Value = $null

# Row 92 (ALC)
$ws.Range("H92").Value = 638.6
$ws.Range("I92").Value = 652.0714
$ws.Range("K92").Value = 652.0714
$ws.Range("M92").Value = 595.9286

# Row 107 (ALC)
$ws.Range("H107").Value = 230
$ws.Range("I107").Value = 400
$ws.Range("J107").Value = 60
$ws.Range("K107").Value = 400
$ws.Range("L107").Value = 60
$ws.Range("M107").Value = 1520
$ws.Range("N107").Value = -3900

# Row 113 (ALC)
$ws.Range("H113").Value = 8000
$ws.Range("J113").Value = 8000
$ws.Range("L113").Value = 8000
$ws.Range("N113").Value = -14508

# Row 116 (ALC)
$ws.Range("H116").Value = 32801
$ws.Range("I116").Value = 12001.667
$ws.Range("J116").Value = 64000
$ws.Range("K116").Value = 12001.667
$ws.Range("L116").Value = 64000
$ws.Range("M116").Value = -8559.666999999999
$ws.Range("N116").Value = -70884

# Row 138 (ALC)
$ws.Range("H138").Value = 4019.375
$ws.Range("I138").Value = 2397.6667
$ws.Range("J138").Value = 4992.4
$ws.Range("K138").Value = 7193.000100000001
$ws.Range("L138").Value = 14977.2
$ws.Range("M138").Value = -2053.000100000001
$ws.Range("N138").Value = -25257.2

$ws = $wb.Worksheets.Item("ARM")
# Row 61 (ARM)
$ws.Range("H61").Value = 5762.4287
$ws.Range("I61").Value = 5472.8335
$ws.Range("J61").Value = 7500
$ws.Range("K61").Value = 5472.8335
$ws.Range("L61").Value = 7500
$ws.Range("M61").Value = -5260.8335
$ws.Range("N61").Value = -7924

# Row 102 (ARM)
$ws.Range("H102").Value = 2896.889
$ws.Range("I102").Value = 2868.5715
$ws.Range("K102").Value = 2868.5715
$ws.Range("M102").Value = -1246.5715

# Row 136 (ARM)
$ws.Range("H136").Value = 5762.4287
$ws.Range("I136").Value = 5472.8335
$ws.Range("J136").Value = 7500
$ws.Range("K136").Value = 16418.5005
$ws.Range("L136").Value = 22500
$ws.Range("M136").Value = -13868.5005
$ws.Range("N136").Value = -27600

$ws = $wb.Worksheets.Item("BSM")
# Row 86 (BSM)
$ws.Range("H86").Value = 9463.333000000001
$ws.Range("J86").Value = 25000
$ws.Range("L86").Value = 25000
$ws.Range("N86").Value = -27246

# Row 89 (BSM)
$ws.Range("H89").Value = 9463.333000000001
$ws.Range("J89").Value = 25000
$ws.Range("L89").Value = 125000
$ws.Range("N89").Value = -136232

# Row 110 (BSM)
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = $null

$ws = $wb.Worksheets.Item("CRP")
# Row 10 (CRP)
$ws.Range("H10").Value = 829.25
$ws.Range("I10").Value = 662
$ws.Range("K10").Value = 662
$ws.Range("M10").Value = -523

# Row 31 (CRP)
$ws.Range("H31").Value = 6000
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 6000
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 6000
$ws.Range("M31").Value = $null
$ws.Range("N31").Value = -6590

# Row 34 (CRP)
$ws.Range("H34").Value = 6000
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 6000
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 6000
$ws.Range("M34").Value = $null
$ws.Range("N34").Value = -6404

# Row 134 (CRP)
$ws.Range("H134").Value = 28000
$ws.Range("I134").Value = 28000
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 84000
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -81465
$ws.Range("N134").Value = $null

$ws = $wb.Worksheets.Item("CUL")
# Row 80 (CUL)
$ws.Range("H80").Value = 14428.571
$ws.Range("J80").Value = 14500
$ws.Range("L80").Value = 43500
$ws.Range("N80").Value = -45372

# Row 83 (CUL)
$ws.Range("H83").Value = 14428.571
$ws.Range("J83").Value = 14500
$ws.Range("L83").Value = 130500
$ws.Range("N83").Value = -139860

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (GSM)
$ws.Range("H70").Value = 6000
$ws.Range("I70").Value = 6000
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -5730
$ws.Range("N70").Value = $null

# Row 73 (GSM)
$ws.Range("H73").Value = 6000
$ws.Range("I73").Value = 6000
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -5064
$ws.Range("N73").Value = $null

$ws = $wb.Worksheets.Item("LTW")
# Row 46 (LTW)
$ws.Range("H46").Value = 9999
$ws.Range("I46").Value = 9999
$ws.Range("K46").Value = 9999
$ws.Range("M46").Value = -9811

# Row 61 (LTW)
$ws.Range("H61").Value = 34001668
$ws.Range("I61").Value = 25502500
$ws.Range("K61").Value = 25502500
$ws.Range("M61").Value = -25502298

# Row 113 (LTW)
$ws.Range("H113").Value = 34001668
$ws.Range("I113").Value = 25502500
$ws.Range("K113").Value = 25502500
$ws.Range("M113").Value = -25500330

# Row 136 (LTW)
$ws.Range("H136").Value = 1699.5
$ws.Range("I136").Value = 1400
$ws.Range("J136").Value = 1999
$ws.Range("K136").Value = 4200
$ws.Range("L136").Value = 5997
$ws.Range("M136").Value = -1650
$ws.Range("N136").Value = -11097

$ws = $wb.Worksheets.Item("WVR")
# Row 21 (WVR)
$ws.Range("H21").Value = 9999
$ws.Range("I21").Value = 9999
$ws.Range("K21").Value = 9999
$ws.Range("M21").Value = -9764

# Row 35 (WVR)
$ws.Range("H35").Value = 9999
$ws.Range("I35").Value = 9999
$ws.Range("K35").Value = 9999
$ws.Range("M35").Value = -9709

# Row 100 (WVR)
$ws.Range("H100").Value = 962.5714
$ws.Range("I100").Value = 897
$ws.Range("K100").Value = 1794
$ws.Range("M100").Value = -1253

# Row 113 (WVR)
$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("K113").Value = 6000
$ws.Range("M113").Value = -3830

# Row 132 (WVR)
$ws.Range("H132").Value = 4249.5
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4249.5
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 12748.5
$ws.Range("M132").Value = $null
$ws.Range("N132").Value = -17808.5

Write-Output "Edits applied successfully"
